# test.xlsx — ProductionPlan sample data refresh
# Renames the Plant/Product id headers, replaces the 5-row numeric sample
# data with a 6-row set that uses text product codes ("x"/"y"), and tweaks
# a couple of view/format bits (column width, selection) to match the
# refreshed workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Plant_Id"
$ws.Range("B1").Value = "Product_Id"

# --- Data rows (2-7) ----------------------------------------------------
# PlantId, ProductId, Target, Unit
$data = @(
    @(31, "x", 20, 1),
    @(32, "x", 11, 1),
    @(33, "x", 13, 1),
    @(34, "x", 14, 1),
    @(35, "x", 50, 1),
    @(11, "y", 60, 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Left-align the Plant/Product/Target columns for the new data rows
$ws.Range("A2:C7").HorizontalAlignment = -4131

# Column B is now text ("x"/"y") rather than numeric codes - widen it to fit
$ws.Columns.Item(2).ColumnWidth = 9.7

# Match the refreshed view's active selection
$ws.Range("B15").Select()
